# TC03_C3DC_phs003164_Race-Unknown.xlsx regression/smoke suite update
# - Simplify the "Treatment Agent" expression in the TreatmentTab query
#   (drop the redundant CONCAT(...) wrapper around REPLACE(...)).
# - Reselect cell B2 (instead of C7) as the active selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$oldText = $treatmentCell.Value2
$newText = $oldText.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$treatmentCell.Value2 = $newText

# Re-apply formatting on the edited cell so it gets its own style entry
# (mirrors the font re-application that happened when the query text was edited).
$treatmentCell.Font.ThemeColor = 1

# Update the sheet's active selection/view.
$ws.Range("B2").Select() | Out-Null
